$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "homecaillotDesktopOutput Algoperfect full testBBA16"

$ws.Range("C2").Value = 99.79513346289541
$ws.Range("D2").Value = 99.64100694689374
$ws.Range("E2").Value = 99.64100694689374
$ws.Range("F2").Value = 99.64100694689374
$ws.Range("G2").Value = 99.64100694689374
$ws.Range("H2").Value = 99.64100694689374
$ws.Range("I2").Value = 99.64100694689374
$ws.Range("J2").Value = 99.64100694689374
$ws.Range("K2").Value = 99.64100694689374
$ws.Range("L2").Value = 99.69663369117117
$ws.Range("M2").Value = 99.69663369117117
$ws.Range("N2").Value = 99.69663369117117
$ws.Range("O2").Value = 99.69663369117117

$wb.Save()
